# Add the new worksheet "Noise Of LIA = -5 dBm" as the last (5th) tab,
# matching the workbook.xml <sheets> order / activeTab in the target diff.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Noise Of LIA = -5 dBm"

# Populate cell data (OSNR/BER/Q-factor table twice - cols A:C and F:H - plus
# the K:L parameter legend), mirroring the "Noise Of LIA = -10 dBm" sheet layout.
$ws.Range("A1").Value = "OSNR, dB"
$ws.Range("B1").Value = "BER"
$ws.Range("C1").Value = "Q-factor"
$ws.Range("F1").Value = "OSNR, dB"
$ws.Range("G1").Value = "BER"
$ws.Range("H1").Value = "Q-factor"
$ws.Range("K1").Value = "Sequence Length, bit"
$ws.Range("L1").Value = 1024
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("K2").Value = "APD M"
$ws.Range("L2").Value = 10
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 13.1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("K3").Value = "TIA Gain, Om"
$ws.Range("L3").Value = 600
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("F4").Value = 13.2
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = "LIA Vpp, V"
$ws.Range("L4").Value = 0.5
$ws.Range("A5").Value = 13
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("F5").Value = 13.3
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("A6").Value = 14
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("F6").Value = 13.4
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("A7").Value = 15
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0
$ws.Range("F7").Value = 13.5
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("A8").Value = 16
$ws.Range("B8").Value = [double]"9.4100000000000005E-40"
$ws.Range("B8").NumberFormat = "0.00E+00"
$ws.Range("C8").Value = 13.14233522
$ws.Range("F8").Value = 13.6
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("A9").Value = 17
$ws.Range("B9").Value = [double]"1.5000000000000001E-33"
$ws.Range("B9").NumberFormat = "0.00E+00"
$ws.Range("C9").Value = 12.011796390000001
$ws.Range("F9").Value = 13.7
$ws.Range("G9").Value = [double]"3.2699999999999999E-18"
$ws.Range("G9").NumberFormat = "0.00E+00"
$ws.Range("H9").Value = 8.6199083549999997
$ws.Range("A10").Value = 18
$ws.Range("B10").Value = [double]"4.34E-38"
$ws.Range("B10").NumberFormat = "0.00E+00"
$ws.Range("C10").Value = 12.84780617
$ws.Range("F10").Value = 13.8
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("A11").Value = 19
$ws.Range("B11").Value = [double]"1.79E-43"
$ws.Range("B11").NumberFormat = "0.00E+00"
$ws.Range("C11").Value = 13.77531804
$ws.Range("F11").Value = 13.9
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("A12").Value = 20
$ws.Range("B12").Value = [double]"6.9099999999999998E-46"
$ws.Range("B12").NumberFormat = "0.00E+00"
$ws.Range("C12").Value = 14.171127439999999
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("A13").Value = 21
$ws.Range("B13").Value = [double]"3.5499999999999997E-46"
$ws.Range("B13").NumberFormat = "0.00E+00"
$ws.Range("C13").Value = 14.217695770000001
$ws.Range("F13").Value = 14.1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("A14").Value = 22
$ws.Range("B14").Value = [double]"8.0999999999999994E-48"
$ws.Range("B14").NumberFormat = "0.00E+00"
$ws.Range("C14").Value = 14.47999132
$ws.Range("F14").Value = 14.2
$ws.Range("G14").Value = [double]"1.4999999999999999E-18"
$ws.Range("G14").NumberFormat = "0.00E+00"
$ws.Range("H14").Value = 8.7082435740000008
$ws.Range("A15").Value = 23
$ws.Range("B15").Value = [double]"2.1900000000000002E-47"
$ws.Range("B15").NumberFormat = "0.00E+00"
$ws.Range("C15").Value = 14.41125551
$ws.Range("F15").Value = 14.3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("A16").Value = 24
$ws.Range("B16").Value = [double]"2.51E-46"
$ws.Range("B16").NumberFormat = "0.00E+00"
$ws.Range("C16").Value = 14.241918310000001
$ws.Range("F16").Value = 14.4
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("A17").Value = 25
$ws.Range("B17").Value = [double]"7.1099999999999997E-47"
$ws.Range("B17").NumberFormat = "0.00E+00"
$ws.Range("C17").Value = 14.32965224
$ws.Range("F17").Value = 14.5
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 0
$ws.Range("A18").Value = 26
$ws.Range("B18").Value = [double]"4.2300000000000003E-48"
$ws.Range("B18").NumberFormat = "0.00E+00"
$ws.Range("C18").Value = 14.52442797
$ws.Range("F18").Value = 14.6
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 0
$ws.Range("A19").Value = 27
$ws.Range("B19").Value = [double]"1.02E-47"
$ws.Range("B19").NumberFormat = "0.00E+00"
$ws.Range("C19").Value = 14.464216070000001
$ws.Range("F19").Value = 14.7
$ws.Range("G19").Value = [double]"2.1399999999999999E-29"
$ws.Range("G19").NumberFormat = "0.00E+00"
$ws.Range("H19").Value = 11.191409569999999
$ws.Range("A20").Value = 28
$ws.Range("B20").Value = [double]"2.73E-49"
$ws.Range("B20").NumberFormat = "0.00E+00"
$ws.Range("C20").Value = 14.71106752
$ws.Range("F20").Value = 14.8
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 0
$ws.Range("A21").Value = 29
$ws.Range("B21").Value = [double]"5.5999999999999999E-47"
$ws.Range("B21").NumberFormat = "0.00E+00"
$ws.Range("C21").Value = 14.346153210000001
$ws.Range("F21").Value = 14.9
$ws.Range("G21").Value = [double]"2.9399999999999999E-28"
$ws.Range("G21").NumberFormat = "0.00E+00"
$ws.Range("H21").Value = 10.959338730000001
$ws.Range("A22").Value = 30
$ws.Range("B22").Value = [double]"1.3100000000000001E-48"
$ws.Range("B22").NumberFormat = "0.00E+00"
$ws.Range("C22").Value = 14.60470864
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("F23").Value = 15.1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("F24").Value = 15.2
$ws.Range("G24").Value = [double]"6.7099999999999997E-29"
$ws.Range("G24").NumberFormat = "0.00E+00"
$ws.Range("H24").Value = 11.08995243
$ws.Range("F25").Value = 15.3
$ws.Range("G25").Value = [double]"2.7399999999999999E-36"
$ws.Range("G25").NumberFormat = "0.00E+00"
$ws.Range("H25").Value = 12.52371668
$ws.Range("F26").Value = 15.4
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("F27").Value = 15.5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("F28").Value = 15.6
$ws.Range("G28").Value = [double]"4.9799999999999997E-25"
$ws.Range("G28").NumberFormat = "0.00E+00"
$ws.Range("H28").Value = 10.264688919999999
$ws.Range("F29").Value = 15.7
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("F30").Value = 15.8
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("F31").Value = 15.9
$ws.Range("G31").Value = [double]"1.03E-32"
$ws.Range("G31").NumberFormat = "0.00E+00"
$ws.Range("H31").Value = 11.852269789999999
$ws.Range("F32").Value = 16
$ws.Range("G32").Value = [double]"9.1100000000000005E-37"
$ws.Range("G32").NumberFormat = "0.00E+00"
$ws.Range("H32").Value = 12.61052355
$ws.Range("F33").Value = 16.100000000000001
$ws.Range("G33").Value = [double]"1.1000000000000001E-25"
$ws.Range("G33").NumberFormat = "0.00E+00"
$ws.Range("H33").Value = 10.40999304
$ws.Range("F34").Value = 16.2
$ws.Range("G34").Value = [double]"3.59E-42"
$ws.Range("G34").NumberFormat = "0.00E+00"
$ws.Range("H34").Value = 13.556331999999999
$ws.Range("F35").Value = 16.3
$ws.Range("G35").Value = [double]"7.3200000000000004E-33"
$ws.Range("G35").NumberFormat = "0.00E+00"
$ws.Range("H35").Value = 11.879482530000001
$ws.Range("F36").Value = 16.399999999999999
$ws.Range("G36").Value = [double]"5.84E-33"
$ws.Range("G36").NumberFormat = "0.00E+00"
$ws.Range("H36").Value = 11.89906201
$ws.Range("F37").Value = 16.5
$ws.Range("G37").Value = [double]"1.6800000000000001E-39"
$ws.Range("G37").NumberFormat = "0.00E+00"
$ws.Range("H37").Value = 13.097228980000001
$ws.Range("F38").Value = 16.600000000000001
$ws.Range("G38").Value = [double]"7.4700000000000001E-39"
$ws.Range("G38").NumberFormat = "0.00E+00"
$ws.Range("H38").Value = 12.982961250000001
$ws.Range("F39").Value = 16.7
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 0
$ws.Range("F40").Value = 16.8
$ws.Range("G40").Value = [double]"4.15E-36"
$ws.Range("G40").NumberFormat = "0.00E+00"
$ws.Range("H40").Value = 12.49012033
$ws.Range("F41").Value = 16.899999999999999
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 0
$ws.Range("F42").Value = 17
$ws.Range("G42").Value = [double]"3.7100000000000001E-33"
$ws.Range("G42").NumberFormat = "0.00E+00"
$ws.Range("H42").Value = 11.938932769999999
$ws.Range("F43").Value = 17.100000000000001
$ws.Range("G43").Value = [double]"3.1999999999999998E-35"
$ws.Range("G43").NumberFormat = "0.00E+00"
$ws.Range("H43").Value = 12.32634011
$ws.Range("F44").Value = 17.2
$ws.Range("G44").Value = [double]"2.02E-39"
$ws.Range("G44").NumberFormat = "0.00E+00"
$ws.Range("H44").Value = 13.08362245
$ws.Range("F45").Value = 17.3
$ws.Range("G45").Value = [double]"1.2900000000000001E-43"
$ws.Range("G45").NumberFormat = "0.00E+00"
$ws.Range("H45").Value = 13.79881484
$ws.Range("F46").Value = 17.399999999999999
$ws.Range("G46").Value = [double]"1.67E-41"
$ws.Range("G46").NumberFormat = "0.00E+00"
$ws.Range("H46").Value = 13.443749950000001
$ws.Range("F47").Value = 17.5
$ws.Range("G47").Value = [double]"8.6800000000000005E-42"
$ws.Range("G47").NumberFormat = "0.00E+00"
$ws.Range("H47").Value = 13.4914346
$ws.Range("F48").Value = 17.600000000000001
$ws.Range("G48").Value = [double]"1.01E-32"
$ws.Range("G48").NumberFormat = "0.00E+00"
$ws.Range("H48").Value = 11.853413359999999
$ws.Range("F49").Value = 17.7
$ws.Range("G49").Value = [double]"2.8399999999999998E-38"
$ws.Range("G49").NumberFormat = "0.00E+00"
$ws.Range("H49").Value = 12.88099411
$ws.Range("F50").Value = 17.8
$ws.Range("G50").Value = [double]"2.1399999999999999E-39"
$ws.Range("G50").NumberFormat = "0.00E+00"
$ws.Range("H50").Value = 13.07901618
$ws.Range("F51").Value = 17.899999999999999
$ws.Range("G51").Value = [double]"4.6799999999999997E-42"
$ws.Range("G51").NumberFormat = "0.00E+00"
$ws.Range("H51").Value = 13.537189229999999
$ws.Range("F52").Value = 18
$ws.Range("G52").Value = [double]"5.7399999999999997E-46"
$ws.Range("G52").NumberFormat = "0.00E+00"
$ws.Range("H52").Value = 14.184003329999999

# Match the saved selection/used-range on the new sheet.
$null = $ws.Range("A1:L52").Select()
